# Idle-time model: fix "Check Count" calculated-column logic so a lone
# (non-begin) idle check immediately after the limit resets to 0 instead of
# being left at 1, and give the "Check Count" table column its formula
# (row 4 previously had a hard-coded 0 with no formula at all).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 4; $r -le 32; $r++) {
    $prev = $r - 1
    $formula = "=IF(AND(Table3[[#This Row],[In Idle Period]],Table3[[#This Row],[Is Begin Idle Period]]), 1, IF(Table3[[#This Row],[In Idle Period]], IF((F$prev+1) <= call_count_limit, IF(AND(F$prev+1=1, NOT(Table3[[#This Row],[Is Begin Idle Period]])), 0, F$prev+1), 0), 0))"
    $ws.Range("F$r").Formula = $formula
}

# Match the selection left behind in the sheet after the edit.
$ws.Range("H13").Select()
